# chore: update Sheets via scheduled runner
#
# Refresh the cached market-board snapshot values (currentAveragePrice,
# currentAveragePriceNQ, currentAveragePriceHQ, LevePriceNQ, LevePriceHQ,
# LeveProfitNQ, LeveProfitHQ -- columns H:N) for the leves whose prices
# moved since the previous run of the scheduled updater.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1317
$ws.Range("I98").Value = 1052.5
$ws.Range("J98").Value = 1493.3334
$ws.Range("K98").Value = 1052.5
$ws.Range("L98").Value = 1493.3334
$ws.Range("M98").Value = 445.5
$ws.Range("N98").Value = -4489.3334

$ws.Range("H105").Value = 30960.25
$ws.Range("J105").Value = 30960.25
$ws.Range("L105").Value = 30960.25
$ws.Range("N105").Value = -37948.25

$ws.Range("H116").Value = 3607.2703
$ws.Range("I116").Value = 2519.15
$ws.Range("J116").Value = 4887.4116
$ws.Range("K116").Value = 2519.15
$ws.Range("L116").Value = 4887.4116
$ws.Range("M116").Value = 922.8499999999999
$ws.Range("N116").Value = -11771.4116

$ws.Range("H122").Value = 1317
$ws.Range("I122").Value = 1052.5
$ws.Range("J122").Value = 1493.3334
$ws.Range("K122").Value = 3157.5
$ws.Range("L122").Value = 4480.0002
$ws.Range("M122").Value = -707.5
$ws.Range("N122").Value = -9380.0002

$ws.Range("H138").Value = 2328.4314
$ws.Range("I138").Value = 1472.9642
$ws.Range("J138").Value = 3369.8696
$ws.Range("K138").Value = 4418.892599999999
$ws.Range("L138").Value = 10109.6088
$ws.Range("M138").Value = 721.1074000000008
$ws.Range("N138").Value = -20389.6088

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4675.5
$ws.Range("I2").Value = 2338.25
$ws.Range("J2").Value = 9350
$ws.Range("K2").Value = 2338.25
$ws.Range("L2").Value = 9350
$ws.Range("M2").Value = -2225.25
$ws.Range("N2").Value = -9576

$ws.Range("H64").Value = 29990
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 29990
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 29990
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -30486

$ws.Range("H67").Value = 29990
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 29990
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 29990
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -31706

$ws.Range("H104").Value = 27499.4
$ws.Range("J104").Value = 27499.4
$ws.Range("L104").Value = 27499.4
$ws.Range("N104").Value = -34487.4

$ws.Range("H106").Value = 44000
$ws.Range("J106").Value = 44000
$ws.Range("L106").Value = 44000
$ws.Range("N106").Value = -46524

$ws.Range("H116").Value = 4675.5
$ws.Range("I116").Value = 2338.25
$ws.Range("J116").Value = 9350
$ws.Range("K116").Value = 2338.25
$ws.Range("L116").Value = 9350
$ws.Range("M116").Value = -44.25
$ws.Range("N116").Value = -13938

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4675.5
$ws.Range("I3").Value = 2338.25
$ws.Range("J3").Value = 9350
$ws.Range("K3").Value = 2338.25
$ws.Range("L3").Value = 9350
$ws.Range("M3").Value = -2224.25
$ws.Range("N3").Value = -9578

$ws.Range("H63").Value = 37235.5
$ws.Range("J63").Value = 37235.5
$ws.Range("L63").Value = 37235.5
$ws.Range("N63").Value = -38607.5

$ws.Range("H66").Value = 37235.5
$ws.Range("J66").Value = 37235.5
$ws.Range("L66").Value = 111706.5
$ws.Range("N66").Value = -118570.5

$ws.Range("H80").Value = 598.4286
$ws.Range("I80").Value = 780.1667
$ws.Range("J80").Value = 462.125
$ws.Range("K80").Value = 780.1667
$ws.Range("L80").Value = 462.125
$ws.Range("M80").Value = 217.8333
$ws.Range("N80").Value = -2458.125

$ws.Range("H83").Value = 598.4286
$ws.Range("I83").Value = 780.1667
$ws.Range("J83").Value = 462.125
$ws.Range("K83").Value = 3900.8335
$ws.Range("L83").Value = 2310.625
$ws.Range("M83").Value = 1091.1665
$ws.Range("N83").Value = -12294.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 8911.714
$ws.Range("I69").Value = 3691
$ws.Range("K69").Value = 3691
$ws.Range("M69").Value = -2942

$ws.Range("H72").Value = 8911.714
$ws.Range("I72").Value = 3691
$ws.Range("K72").Value = 11073
$ws.Range("M72").Value = -7329

$ws.Range("H93").Value = 23650
$ws.Range("I93").Value = 17500
$ws.Range("J93").Value = 29800
$ws.Range("K93").Value = 17500
$ws.Range("L93").Value = 29800
$ws.Range("M93").Value = -15628
$ws.Range("N93").Value = -33544

$ws.Range("H103").Value = 31944
$ws.Range("I103").Value = 24000
$ws.Range("J103").Value = 39888
$ws.Range("K103").Value = 24000
$ws.Range("L103").Value = 39888
$ws.Range("M103").Value = -22828
$ws.Range("N103").Value = -42232

$ws.Range("H134").Value = 12874
$ws.Range("I134").Value = 9194.666999999999
$ws.Range("J134").Value = 51507
$ws.Range("K134").Value = 27584.001
$ws.Range("L134").Value = 154521
$ws.Range("M134").Value = -25049.001
$ws.Range("N134").Value = -159591

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 1346.4445
$ws.Range("I45").Value = 500
$ws.Range("J45").Value = 1452.25
$ws.Range("K45").Value = 1500
$ws.Range("L45").Value = 4356.75
$ws.Range("M45").Value = -968
$ws.Range("N45").Value = -5420.75

$ws.Range("H68").Value = 1673.4546
$ws.Range("I68").Value = 1413
$ws.Range("J68").Value = 1731.3334
$ws.Range("K68").Value = 4239
$ws.Range("L68").Value = 5194.0002
$ws.Range("M68").Value = -3428
$ws.Range("N68").Value = -6816.0002

$ws.Range("H71").Value = 1673.4546
$ws.Range("I71").Value = 1413
$ws.Range("J71").Value = 1731.3334
$ws.Range("K71").Value = 12717
$ws.Range("L71").Value = 15582.0006
$ws.Range("M71").Value = -8661
$ws.Range("N71").Value = -23694.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 574.2692
$ws.Range("I107").Value = 372.66666
$ws.Range("J107").Value = 849.1818
$ws.Range("K107").Value = 372.66666
$ws.Range("L107").Value = 849.1818
$ws.Range("M107").Value = 1547.33334
$ws.Range("N107").Value = -4689.1818

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 33049.6
$ws.Range("J62").Value = 33049.6
$ws.Range("L62").Value = 33049.6
$ws.Range("N62").Value = -34297.6

$ws.Range("H65").Value = 33049.6
$ws.Range("J65").Value = 33049.6
$ws.Range("L65").Value = 99148.79999999999
$ws.Range("N65").Value = -105388.8

$ws.Range("H105").Value = 44000
$ws.Range("J105").Value = 44000
$ws.Range("L105").Value = 44000
$ws.Range("N105").Value = -50988

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 4600

$ws.Range("H104").Value = 40000
$ws.Range("J104").Value = 40000
$ws.Range("L104").Value = 40000
$ws.Range("N104").Value = -46988

$ws.Range("H105").Value = 44990.75
$ws.Range("J105").Value = 44990.75
$ws.Range("L105").Value = 44990.75
$ws.Range("N105").Value = -51978.75
